$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tariff Added: append a new data row (row 3) below the existing one,
# following the same 25-column layout (meterId + previous peak/off-peak values).
$ws.Cells.Item(3, 1).Value = 12345678901122
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 0
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 1798
$ws.Cells.Item(3, 19).Value = 9325
$ws.Cells.Item(3, 20).Value = 1800
$ws.Cells.Item(3, 21).Value = 9332
$ws.Cells.Item(3, 22).Value = 0
$ws.Cells.Item(3, 23).Value = 0
$ws.Cells.Item(3, 24).Value = 1849
$ws.Cells.Item(3, 25).Value = 9571

# Searchbox updated: scroll the view over to column Q and move the
# active selection to U4.
$win = $excel.ActiveWindow
$win.ScrollColumn = 17
$win.ScrollRow = 1
$ws.Range("U4").Select()
